$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data block (Plátano, Vega Monumental Concepción) is added.
# It belongs right above the current "44421" block (row 356), so insert a
# fresh row there - this shifts every following row down by one and the
# last existing row (367) ends up duplicated into the newly created row 368.
$ws.Rows.Item(356).Insert()

$ws.Cells.Item(356, 1).Value  = 11
$ws.Cells.Item(356, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(356, 3).Value  = "Bíobío"
$ws.Cells.Item(356, 4).Value  = 44509
$ws.Cells.Item(356, 5).Value  = 8
$ws.Cells.Item(356, 6).Value  = "Fruta"
$ws.Cells.Item(356, 7).Value  = 100108
$ws.Cells.Item(356, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(356, 9).Value  = 100108006
$ws.Cells.Item(356, 10).Value = "Plátano"
$ws.Cells.Item(356, 11).Value = "Sin especificar"
$ws.Cells.Item(356, 12).Value = "Primera Pintón"
$ws.Cells.Item(356, 13).Value = 1100
$ws.Cells.Item(356, 14).Value = 17000
$ws.Cells.Item(356, 15).Value = 18000
$ws.Cells.Item(356, 16).Value = 17545
$ws.Cells.Item(356, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(356, 18).Value = "Ecuador"
$ws.Cells.Item(356, 19).Value = 877
$ws.Cells.Item(356, 20).Value = 20
